$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 509 (existing rows 509:519 shift down to 512:522).
$ws.Rows.Item(509).Insert()
$ws.Rows.Item(509).Insert()
$ws.Rows.Item(509).Insert()

# Columns A, B, C, E, F, G, H, I, J, K are constant for this whole block of
# "Arandano (blue)" rows at Mercado Mayorista Lo Valledor de Santiago.
$constA = 6
$constB = "Mercado Mayorista Lo Valledor de Santiago"
$constC = "Metropolitana"
$constE = 13
$constF = "Fruta"
$constG = 100101
$constH = "Berries"
$constI = 100101001
$constJ = "Arándano (blue)"
$constK = "Sin especificar"

function Fill-Row($r, $d, $l, $m, $n, $o, $p, $q, $rr, $s, $t) {
    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $constI
    $ws.Cells.Item($r, 10).Value = $constJ
    $ws.Cells.Item($r, 11).Value = $constK
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rr
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}

# New row 509
Fill-Row 509 44939 "Especial" 300 3000 3000 3000 "`$/bandeja 2 kilos" "Provincia de Curicó" 1500 2

# New row 510
Fill-Row 510 44939 "Especial" 2000 3000 3000 3000 "`$/bandeja 2 kilos" "Región del Maule" 1500 2

# New row 511
Fill-Row 511 44939 "Segunda" 250 2000 2000 2000 "`$/bandeja 2 kilos" "Provincia de Curicó" 1000 2
